$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 3 (market size inputs)
$ws.Range("B3").Formula = "=1.5"
$ws.Range("C3").Value = 1.2
$ws.Range("D3").Value = 1.1
$ws.Range("E3").Value = 1.05
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6
$ws.Range("H3").Value = 0.5

# Add style to new empty rows 6-8 (use same style as row3 col A, which is s="1")
$ws.Range("A6:F6").Style = $ws.Range("A3").Style
$ws.Range("A7:H7").Style = $ws.Range("A3").Style
$ws.Range("A8:F8").Style = $ws.Range("A3").Style
